$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("comforter-cda")

# Fill in the previously-missing Start Time / End Time values for row 92
$ws.Range("B92").Value = 0
$ws.Range("C92").Value = 0

# Add new daily power record row 93 - copy formatting from row 92's
# populated columns first (A, D, E, F) so the new row matches the
# table's existing look; B93/C93 are intentionally left blank, same
# as row 92 was before it got its Start/End Time values.
$ws.Range("A92").Copy()
$ws.Range("A93").PasteSpecial(-4122)
$ws.Range("D92:F92").Copy()
$ws.Range("D93:F93").PasteSpecial(-4122)

$ws.Range("A93").Value = 43417
$ws.Range("D93").Formula = "=(C93-B93)* 1440"
$ws.Range("E93").Formula = "=IF(C93>B93, (C93-B93)*1440, (B93-C93)*1440)"
$ws.Range("F93").Formula = "=ABS((C93-B93)*1440)"

# Resize the table to include the new row
$table = $ws.ListObjects.Item("comforter_cda_table")
$table.Resize($ws.Range("A1:F93"))

# Update the selection to reflect the newly added row, matching Excel's
# auto-scroll/selection behavior when data is appended at the bottom
$ws.Range("B93").Select()
